$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") "26.782.63"
Set-TextValue $ws.Range("E2") "  +0.18%  "
Set-TextValue $ws.Range("D3") "1.643.94"
Set-TextValue $ws.Range("E3") "  -0.26%  "
Set-TextValue $ws.Range("E4") "  +0.68%  "
Set-TextValue $ws.Range("D5") "216.84"
Set-TextValue $ws.Range("E5") "  +0.35%  "
Set-TextValue $ws.Range("E6") "  -0.37%  "
Set-TextValue $ws.Range("E7") "  +0.58%  "
Set-TextValue $ws.Range("E8") "  -0.82%  "
Set-TextValue $ws.Range("D9") "0.0627"
Set-TextValue $ws.Range("E9") "  -0.27%  "
Set-TextValue $ws.Range("D10") "19.18"
Set-TextValue $ws.Range("E10") "  -1.04%  "
Set-TextValue $ws.Range("D11") "0.0842"
Set-TextValue $ws.Range("E11") "  -0.35%  "
Set-TextValue $ws.Range("D12") "1.869.22"
Set-TextValue $ws.Range("D13") "1.635.96"
Set-TextValue $ws.Range("E13") "  -0.60%  "
Set-TextValue $ws.Range("E14") "  -1.13%  "
Set-TextValue $ws.Range("D15") "0.527"
Set-TextValue $ws.Range("E15") "  -1.40%  "
Set-TextValue $ws.Range("D16") "64.53"
Set-TextValue $ws.Range("E16") "  -2.86%  "
Set-TextValue $ws.Range("D17") "26.793.72"
Set-TextValue $ws.Range("E17") "  +0.19%  "
Set-TextValue $ws.Range("E18") "  -2.35%  "
Set-TextValue $ws.Range("D19") "214.25"
Set-TextValue $ws.Range("E19") "  -2.77%  "
Set-TextValue $ws.Range("E20") "  +0.64%  "
Set-TextValue $ws.Range("E21") "  -0.61%  "
Set-TextValue $ws.Range("D22") "2.40"
Set-TextValue $ws.Range("E22") "  +13.00%  "
Set-TextValue $ws.Range("D23") "6.28"
Set-TextValue $ws.Range("E23") "  -0.97%  "
Set-TextValue $ws.Range("E24") "  -2.14%  "
Set-TextValue $ws.Range("D25") "145.03"
Set-TextValue $ws.Range("E25") "  -1.50%  "
Set-TextValue $ws.Range("E26") "  +0.96%  "
Set-TextValue $ws.Range("E27") "  -2.04%  "
Set-TextValue $ws.Range("E28") "  -0.13%  "
Set-TextValue $ws.Range("E29") "  -1.38%  "
Set-TextValue $ws.Range("E30") "  -1.12%  "
Set-TextValue $ws.Range("E31") "  +0.17%  "
Set-TextValue $ws.Range("E32") "  -3.05%  "
Set-TextValue $ws.Range("E33") "  -1.74%  "
Set-TextValue $ws.Range("D34") "1.289.50"
Set-TextValue $ws.Range("E34") "  +0.00%  "
Set-TextValue $ws.Range("E35") "  -1.08%  "
Set-TextValue $ws.Range("E36") "  +1.07%  "
Set-TextValue $ws.Range("D37") "0.0176"
Set-TextValue $ws.Range("E37") "  -4.95%  "
Set-TextValue $ws.Range("D38") "0.537"
Set-TextValue $ws.Range("E38") "  +2.50%  "
Set-TextValue $ws.Range("D39") "0.826"
Set-TextValue $ws.Range("E39") "  -0.51%  "
Set-TextValue $ws.Range("E40") "  +0.46%  "
Set-TextValue $ws.Range("D41") "0.808"
Set-TextValue $ws.Range("E41") "  -0.32%  "
Set-TextValue $ws.Range("E42") "  -0.24%  "
Set-TextValue $ws.Range("D43") "5.36"
Set-TextValue $ws.Range("E43") "  -1.30%  "
Set-TextValue $ws.Range("D44") "1.795.63"
Set-TextValue $ws.Range("E44") "  +0.43%  "
Set-TextValue $ws.Range("B45") "Aave"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "60.44"
Set-TextValue $ws.Range("E45") "  +1.66%  "
Set-TextValue $ws.Range("B46") "Quant"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D46") "91.24"
Set-TextValue $ws.Range("E46") "  -2.67%  "
Set-TextValue $ws.Range("E47") "  -0.26%  "
Set-TextValue $ws.Range("B48") "BabyDogeCoin"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D48") "0.0₆0104"
Set-TextValue $ws.Range("E48") "  -1.84%  "
Set-TextValue $ws.Range("B49") "Cronos"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D49") "0.0520"
Set-TextValue $ws.Range("E49") "  +0.66%  "
Set-TextValue $ws.Range("B50") "EnergySwap"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.70"
Set-TextValue $ws.Range("E50") "  -1.21%  "
Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.0979"
Set-TextValue $ws.Range("E51") "  +0.13%  "
